# Remove the leading "▲" run that precedes the "表 8-2-6 成就徽章"
# caption paragraph, leaving the caption run itself untouched.
#
# The document contains several "▲表 8-2-x ..." captions (8-2-6 through
# 8-2-10); only the very first one (8-2-6, "成就徽章") loses its "▲" run,
# so we anchor the edit on that caption's exact, unique text rather than
# a fragile paragraph index.

$d = $word.ActiveDocument
$target = $d.Content
$found = $target.Find.Execute("▲表 8-2-6 成就徽章", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $target now spans the whole matched caption; its first character is
    # the standalone "▲" run - delete just that character so the
    # following run ("表 8-2-6 成就徽章") is left completely intact.
    $target.Characters(1).Delete()
}
